$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 (2013年) values
$ws.Range("B5").Value = 8517893
$ws.Range("C5").Value = 8314230
$ws.Range("D5").Value = 117150
$ws.Range("E5").Value = 106387
$ws.Range("F5").Value = 1510930
$ws.Range("G5").Value = 97276
$ws.Range("H5").Value = 5738666
$ws.Range("I5").Value = 20106
$ws.Range("J5").Value = 65950
$ws.Range("K5").Value = 125325
$ws.Range("L5").Value = 134013

# Add new row 12 (2021年)
$ws.Range("A11:L11").Copy()
$ws.Range("A12:L12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 28665212
$ws.Range("C12").Value = 28383434
$ws.Range("D12").Value = 78357
$ws.Range("E12").Value = 135203
$ws.Range("F12").Value = 1711337
$ws.Range("G12").Value = 146575
$ws.Range("H12").Value = 26288321
$ws.Range("I12").Value = 5631
$ws.Range("J12").Value = 31729
$ws.Range("K12").Value = 106993
$ws.Range("L12").Value = 103534
